# Additional companies sent for questionaire
# Remove the "Parent company" (column B) and "Location County/City"
# (column E) columns from the locomotive list, shifting remaining
# columns to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from right to left so column letters of not-yet-deleted
# columns remain stable.
$ws.Range("E1").EntireColumn.Delete()
$ws.Range("B1").EntireColumn.Delete()

# Restore the selection/active cell as captured in the saved workbook.
$ws.Range("I8").Select()
